# [AFG] added final excel sheets for Afghanistan
#
# Adds a new worksheet "ODI Batting Extra" (3rd sheet, after "Player Info"
# and "ODI Batting") containing extended ODI batting stats for player 4615
# (Hazratullah Zazai): MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6,
# PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH.

$wb = $excel.ActiveWorkbook

# --- Create the new worksheet at the end of the workbook ---------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# --- Header row (A1:F1) --------------------------------------------------
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $headers.Length; $col++) {
    $newSheet.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Reuse the exact same header style (bold, centered, thin border) already
# used by the "ID"/"NAME"/... headers on the other sheets.
$wb.Worksheets.Item("Player Info").Range("A1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows -------------------------------------------------------------
# Columns: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
# BATTING_POSITION is numeric when known; NUM_4/NUM_6/PERCENT_RUNS_OF_TOTAL are
# stored as text (mirrors the source data, which keeps them as strings even
# though several look numeric), and may be blank for a few matches.
$rows = @(
    @("4189", "1", "1", "1", "6.17%"),
    @("4190", "2", "0", "0", "0.55%"),
    @("4257", "2", "3", "1", "15.15%"),
    @("4259", "",  "",  "",  ""),
    @("4262", "",  "",  "",  ""),
    @("4265", "1", "1", "1", "6.28%"),
    @("4290", "2", "2", "1", "5.20%"),
    @("4299", "2", "2", "0", "10.14%"),
    @("4306", "2", "0", "0", ""),
    @("4309", "2", "3", "1", "19.74%"),
    @("4315", "1", "5", "1", "19.77%"),
    @("4323", "1", "3", "0", "17.60%"),
    @("4332", "",  "",  "",  ""),
    @("4377", "1", "1", "0", "4.64%"),
    @("4378", "1", "1", "2", "11.50%"),
    @("4379", "1", "7", "2", "20.08%")
)

$rowIndex = 2
foreach ($r in $rows) {
    $matchCode = $r[0]
    $battingPosition = $r[1]
    $num4 = $r[2]
    $num6 = $r[3]
    $percentRuns = $r[4]

    # MATCH_CODE (column A) - always text, even though it looks numeric
    $cA = $newSheet.Cells.Item($rowIndex, 1)
    $cA.NumberFormat = "@"
    $cA.Value = $matchCode
    $cA.Style = "Normal"

    # BATTING_POSITION (column B) - numeric when present, blank otherwise
    $cB = $newSheet.Cells.Item($rowIndex, 2)
    if ($battingPosition -eq "") {
        $cB.NumberFormat = "@"
        $cB.Value = ""
        $cB.Style = "Normal"
    } else {
        $cB.Value = [int]$battingPosition
    }

    # NUM_4 (column C) - text
    $cC = $newSheet.Cells.Item($rowIndex, 3)
    $cC.NumberFormat = "@"
    $cC.Value = $num4
    $cC.Style = "Normal"

    # NUM_6 (column D) - text
    $cD = $newSheet.Cells.Item($rowIndex, 4)
    $cD.NumberFormat = "@"
    $cD.Value = $num6
    $cD.Style = "Normal"

    # PERCENT_RUNS_OF_TOTAL (column E) - text
    $cE = $newSheet.Cells.Item($rowIndex, 5)
    $cE.NumberFormat = "@"
    $cE.Value = $percentRuns
    $cE.Style = "Normal"

    # MAN_OF_MATCH (column F) - always "NO" for this player
    $newSheet.Cells.Item($rowIndex, 6).Value = "NO"

    $rowIndex++
}

# Restore the original active sheet/selection state ("Player Info" was the
# active sheet before this edit).
$wb.Worksheets.Item("Player Info").Activate()
